# Updating the model to V7 - reducing the hidden layers.
# Adds a "hidden layers" annotation column (K) for the existing model rows,
# appends a new model_V7 summary row (10), and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New annotation column K for existing rows -----------------------------
# Assignment order matters for shared-string allocation, so write these in
# the same order the new unique strings first appear in the target file:
# "Lower hidden layers", "model_V7", "4 hidden layers", "6 hidden layer",
# "8 hidden layers".
$ws.Range("K10").Value = "Lower hidden layers"
$ws.Range("A10").Value = "model_V7"
$ws.Range("K6").Value = "4 hidden layers"
$ws.Range("K8").Value = "6 hidden layer"
$ws.Range("K7").Value = "8 hidden layers"
$ws.Range("K9").Value = "4 hidden layers"

# --- New row 10: model_V7 summary data --------------------------------------
$ws.Range("B10").Value = 382
$ws.Range("C10").Value = 625
$ws.Range("D10").Value = 607
$ws.Range("E10").Value = 386
$ws.Range("F10").Formula = "=SUM(B10:E10)"
$ws.Range("G10").Formula = "=E10/F10"
$ws.Range("I10").Formula = "=E10+D10"
$ws.Range("J10").Formula = "=I10/F10"

# --- Update the current selection/view --------------------------------------
$ws.Range("E7:E10").Select() | Out-Null
